$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price column (some values look numeric
# and would otherwise be auto-converted by Excel on assignment).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 40 and 41 swap places (Stacks moves up, FirstDigitalUSD moves down)
# with updated price/volume figures.
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "2.24"
$ws.Range("E40").Value = "  -3.11%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.14%  "

# Updated price / 1h volume figures for the remaining rows.
$ws.Range("D2").Value = "69.343.15"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "3.687.80"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "677.94"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("D11").Value = "0.441"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "4.307.63"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "32.36"
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "3.678.91"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "69.313.45"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("D18").Value = "16.00"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "468.68"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "3.833.81"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -4.71%  "
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("D34").Value = "26.92"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "3.677.37"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "6.31"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").Value = "170.48"
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("D44").Value = "0.943"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").Value = "47.17"
$ws.Range("E45").Value = "  -2.34%  "
$ws.Range("E46").Value = "  -6.81%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").Value = "2.70"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").Value = "0.000277"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("E51").Value = "  -2.76%  "
